$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-14 03:46:47"
$wsZh.Range("H2").Value = "2016-03-14 03:47:10"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-14 03:46:50"
$wsDe.Range("H2").Value = "2016-03-14 03:47:16"
